# Auto-generated edit script applying the diff to Bahamut_Profits sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 265833.1
$ws.Range("J17").Value = 265833.1
$ws.Range("L17").Value = 797499.2999999999
$ws.Range("N17").Value = -797835.2999999999
$ws.Range("H98").Value = 1248
$ws.Range("I98").Value = 1367.0714
$ws.Range("J98").Value = 581.2
$ws.Range("K98").Value = 1367.0714
$ws.Range("L98").Value = 581.2
$ws.Range("M98").Value = 130.9286
$ws.Range("N98").Value = -3577.2
$ws.Range("H122").Value = 1248
$ws.Range("I122").Value = 1367.0714
$ws.Range("J122").Value = 581.2
$ws.Range("K122").Value = 4101.2142
$ws.Range("L122").Value = 1743.6
$ws.Range("M122").Value = -1651.2142
$ws.Range("N122").Value = -6643.6
$ws.Range("H129").Value = 2315858.8
$ws.Range("I129").Value = 395
$ws.Range("J129").Value = 2470223
$ws.Range("K129").Value = 1185
$ws.Range("L129").Value = 7410669
$ws.Range("M129").Value = 3815
$ws.Range("N129").Value = -7420669
$ws.Range("H137").Value = 1073.3334
$ws.Range("I137").Value = 927.2727
$ws.Range("J137").Value = 1475
$ws.Range("K137").Value = 2781.8181
$ws.Range("L137").Value = 4425
$ws.Range("M137").Value = -231.8181
$ws.Range("N137").Value = -9525
$ws.Range("H138").Value = 3738.23
$ws.Range("J138").Value = 4476.9375
$ws.Range("L138").Value = 13430.8125
$ws.Range("N138").Value = -23710.8125
$ws.Range("H139").Value = 48454
$ws.Range("J139").Value = 48454
$ws.Range("L139").Value = 48454
$ws.Range("N139").Value = -58734

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H139").Value = 45695
$ws.Range("J139").Value = 45695
$ws.Range("L139").Value = 45695
$ws.Range("N139").Value = -55975

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1275.8823
$ws.Range("I94").Value = 1060
$ws.Range("J94").Value = 1977.5
$ws.Range("K94").Value = 1060
$ws.Range("L94").Value = 1977.5
$ws.Range("M94").Value = -609
$ws.Range("N94").Value = -2879.5
$ws.Range("H105").Value = 4822.6313
$ws.Range("I105").Value = 4803.3335
$ws.Range("J105").Value = 4895
$ws.Range("K105").Value = 4803.3335
$ws.Range("L105").Value = 4895
$ws.Range("M105").Value = -3056.3335
$ws.Range("N105").Value = -8389
$ws.Range("H135").Value = 50319.75
$ws.Range("J135").Value = 50319.75
$ws.Range("L135").Value = 50319.75
$ws.Range("N135").Value = -60459.75
$ws.Range("H138").Value = 52780
$ws.Range("J138").Value = 52780
$ws.Range("L138").Value = 52780
$ws.Range("N138").Value = -63060

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4349
$ws.Range("I31").Value = 5132
$ws.Range("J31").Value = 2000
$ws.Range("K31").Value = 5132
$ws.Range("L31").Value = 2000
$ws.Range("M31").Value = -4837
$ws.Range("N31").Value = -2590
$ws.Range("H34").Value = 4349
$ws.Range("I34").Value = 5132
$ws.Range("J34").Value = 2000
$ws.Range("K34").Value = 5132
$ws.Range("L34").Value = 2000
$ws.Range("M34").Value = -4930
$ws.Range("N34").Value = -2404
$ws.Range("H122").Value = 559.5
$ws.Range("I122").Value = 588.8889
$ws.Range("J122").Value = 471.33334
$ws.Range("K122").Value = 1766.6667
$ws.Range("L122").Value = 1414.00002
$ws.Range("M122").Value = 683.3332999999998
$ws.Range("N122").Value = -6314.000019999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1398
$ws.Range("I5").Value = 1094.3334
$ws.Range("J5").Value = 2569.2856
$ws.Range("K5").Value = 3283.0002
$ws.Range("L5").Value = 7707.8568
$ws.Range("M5").Value = -3171.0002
$ws.Range("N5").Value = -7931.8568
$ws.Range("H101").Value = 4812.5
$ws.Range("J101").Value = 5357.143
$ws.Range("L101").Value = 16071.429
$ws.Range("N101").Value = -20939.429
$ws.Range("H122").Value = 358112.56
$ws.Range("I122").Value = 645.2308
$ws.Range("J122").Value = 667917.6
$ws.Range("K122").Value = 5807.077200000001
$ws.Range("L122").Value = 6011258.399999999
$ws.Range("M122").Value = -3357.077200000001
$ws.Range("N122").Value = -6016158.399999999
$ws.Range("H135").Value = 1398
$ws.Range("I135").Value = 1094.3334
$ws.Range("J135").Value = 2569.2856
$ws.Range("K135").Value = 9849.000599999999
$ws.Range("L135").Value = 23123.5704
$ws.Range("M135").Value = -7314.000599999999
$ws.Range("N135").Value = -28193.5704

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 10000050
$ws.Range("I10").Value = 100
$ws.Range("J10").Value = 20000000
$ws.Range("K10").Value = 100
$ws.Range("L10").Value = 20000000
$ws.Range("M10").Value = 69
$ws.Range("N10").Value = -20000338
$ws.Range("H70").Value = 5658.2856
$ws.Range("I70").Value = 4321.6
$ws.Range("J70").Value = 9000
$ws.Range("K70").Value = 4321.6
$ws.Range("L70").Value = 9000
$ws.Range("M70").Value = -4051.6
$ws.Range("N70").Value = -9540
$ws.Range("H73").Value = 5658.2856
$ws.Range("I73").Value = 4321.6
$ws.Range("J73").Value = 9000
$ws.Range("K73").Value = 4321.6
$ws.Range("L73").Value = 9000
$ws.Range("M73").Value = -3385.6
$ws.Range("N73").Value = -10872
$ws.Range("H136").Value = 21830.6
$ws.Range("J136").Value = 21830.6
$ws.Range("L136").Value = 65491.8
$ws.Range("N136").Value = -70591.79999999999
$ws.Range("H138").Value = 41833.332
$ws.Range("J138").Value = 45250
$ws.Range("L138").Value = 45250
$ws.Range("N138").Value = -55530

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 900
$ws.Range("I58").Value = 666.6667
$ws.Range("K58").Value = 666.6667
$ws.Range("M58").Value = -406.6667
$ws.Range("H132").Value = 2947.275
$ws.Range("I132").Value = 2705.1035
$ws.Range("J132").Value = 3585.7273
$ws.Range("K132").Value = 8115.310500000001
$ws.Range("L132").Value = 10757.1819
$ws.Range("M132").Value = -5585.310500000001
$ws.Range("N132").Value = -15817.1819
$ws.Range("H138").Value = 38026.668
$ws.Range("J138").Value = 38026.668
$ws.Range("L138").Value = 38026.668
$ws.Range("N138").Value = -48306.668

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1400
$ws.Range("I136").Value = 1343.1082
$ws.Range("J136").Value = 1926.25
$ws.Range("K136").Value = 4029.3246
$ws.Range("L136").Value = 5778.75
$ws.Range("M136").Value = -1479.3246
$ws.Range("N136").Value = -10878.75
$ws.Range("H138").Value = 49583.332
$ws.Range("J138").Value = 49583.332
$ws.Range("L138").Value = 49583.332
$ws.Range("N138").Value = -59863.332
